# Applies the content-rotation edit described by the LOQ4063.docx diff.
# Most of the change is several whole-paragraph bodies swapping places; the
# 'Avaliacao' paragraph additionally has three value-runs (Metodo / Criterio /
# Norma de recuperacao) updated in place so the bold labels and each run's
# trailing <w:br/> line break stay untouched.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 'Avaliacao' paragraph: swap the three labelled value runs in place.
# Because 'Metodo' ends up holding the *old* 'Criterio' text and
# 'Criterio' ends up holding the *old* 'Norma de recuperacao' text,
# a naive repeated Find/Replace would match its own freshly-written
# output on the next call. Instead, locate the three value spans by
# their fixed labels up front, then overwrite the spans back-to-front
# (right-most span first) via Range(start, end) so earlier offsets
# never shift under us.
# ------------------------------------------------------------------

$avaliacaoPara = $d.Paragraphs.Item(14)
$paraStart = $avaliacaoPara.Range.Start
$fullText = $avaliacaoPara.Range.Text

$metodoLabel = 'Método: '
$criterioLabel = 'Critério: '
$normaLabel = 'Norma de recuperação: '

$iMetodo = $fullText.IndexOf($metodoLabel)
$iCriterio = $fullText.IndexOf($criterioLabel)
$iNorma = $fullText.IndexOf($normaLabel)

$metodoValStart = $iMetodo + $metodoLabel.Length
$metodoValEnd = $fullText.IndexOf([char]11, $metodoValStart)
$criterioValStart = $iCriterio + $criterioLabel.Length
$criterioValEnd = $fullText.IndexOf([char]11, $criterioValStart)
$normaValStart = $iNorma + $normaLabel.Length
$normaValEnd = $fullText.Length

$normaNewValue = '1) FOUST, Alan S.; WENZEL, Leonard A.; CLUMP, Curtis W.; MAUS, Louis; ANDERSEN, L. Bryce. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 1982.' + [char]11 + '2) GEANKOPLIS, Christie John. Transport Processes and Separation Process Principles. New York: Prentice Hall, 2003.' + [char]11 + '3) COUPER, James R.; PENNEY, W. Roy; FAIR, James R.; WALAS, Stanley M. Chemical Process Equipment: Selection and Design. Amsterdam: Elsevier, 2005.' + [char]11 + '4) FOGLER, H. S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2002.' + [char]11 + '5) LEVENSPIEL, O. Chemical Reaction Engineering. 3rd.ed. New York: John Wiley & Sons, 1998.' + [char]11 + '6) PERRY, Robert H.; GREEN, Don W. Perry''s Chemical Engineers'' Handbook. 8th.ed. New York: McGraw-Hill, 2008.'
$criterioNewValue = 'A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação'
$metodoNewValue = 'A média do período será definida pelo professor da disciplina. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental).'

$d.Range($paraStart + $normaValStart, $paraStart + $normaValEnd).Text = $normaNewValue
$d.Range($paraStart + $criterioValStart, $paraStart + $criterioValEnd).Text = $criterioNewValue
$d.Range($paraStart + $metodoValStart, $paraStart + $metodoValEnd).Text = $metodoNewValue

# ------------------------------------------------------------------
# Remaining whole-paragraph-run content rotation. Each of these runs
# is a standalone paragraph (no sibling runs / bold labels), so a
# plain Range.Text assignment is safe and keeps every other
# paragraph's index stable (paragraph count never changes).
# A vertical-tab char (decimal 11) is Word's in-memory marker for a
# manual line break and round-trips to <w:br/> on save.
# ------------------------------------------------------------------

# 'Objetivos' body paragraph -> short 3-item 'Programa resumido' style list
$d.Paragraphs.Item(6).Range.Text = '1) Reatores químicos' + [char]11 + '2) Fermentação' + [char]11 + '3) Processos químicos'

# 'Docente(s) Responsavel(eis)' paragraph -> long 'Objetivos' paragraph text
$d.Paragraphs.Item(8).Range.Text = 'Experiências em laboratório de caráter multidisciplinar que tem por objetivo colocar o aluno em contato com equipamentos de engenharia e consolidar os conceitos de fenômenos de transporte. O desenvolvimento das atividades inclui montagem, medidas e interpretação de resultados em áreas relevantes da engenharia como cinética e reatores químicos, fenômenos de transporte, operações unitárias e processos químicos industriais. A disciplina permite um programa dinâmico, onde os experimentos poderão ser mudados e/ou revezados em função da evolução dos laboratórios ou necessidades específicas.'

# 'Programa resumido' paragraph -> long 3-item 'Programa' detail list
$d.Paragraphs.Item(10).Range.Text = '1) Reatores químicos: operação de reator de mistura para avaliar a influência do tempo de residência na conversão.' + [char]11 + '2) Fermentação: determinação de parâmetros cinéticos em processo de fermentação alcoólica por leveduras em reatores bioquímicos.' + [char]11 + '3) Processos químicos: executar síntese de produto de interesse industrial em uma planta de pequeno porte. Poderão ser observados aspectos como instrumentação (controles de vazão, temperatura, nível, etc), equipamentos diversos de operações unitárias, sistema de aquisição de dados, etc'

# 'Programa' paragraph -> old 'Metodo:' value text
$d.Paragraphs.Item(12).Range.Text = 'Aplicação de prova(s) e relatório(s).'

# 'Bibliografia' paragraph -> docente identification line
$d.Paragraphs.Item(16).Range.Text = '5816812 - João Paulo Alves Silva'

